# Update the "addition and subtraction within 100" answer table:
# replace each old equation with its new equation, one occurrence at a
# time (MatchWholeWord, Replace:=wdReplaceOne) so that cells sharing a
# resulting value are never re-matched after being set.
$d = $word.ActiveDocument

$d.Content.Find.Execute("24+38=62", $true, $false, $false, $false, $false, $true, 1, $false, "94-22=72", 2) | Out-Null
$d.Content.Find.Execute("94+2=96", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("16+2=18", $true, $false, $false, $false, $false, $true, 1, $false, "25-8=17", 2) | Out-Null
$d.Content.Find.Execute("76-22=54", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=79", 2) | Out-Null
$d.Content.Find.Execute("12+74=86", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=71", 2) | Out-Null
$d.Content.Find.Execute("84-81=3", $true, $false, $false, $false, $false, $true, 1, $false, "23-0=23", 2) | Out-Null
$d.Content.Find.Execute("77-34=43", $true, $false, $false, $false, $false, $true, 1, $false, "75-30=45", 2) | Out-Null
$d.Content.Find.Execute("54-9=45", $true, $false, $false, $false, $false, $true, 1, $false, "35-34=1", 2) | Out-Null
$d.Content.Find.Execute("43+3=46", $true, $false, $false, $false, $false, $true, 1, $false, "77-9=68", 2) | Out-Null
$d.Content.Find.Execute("36-28=8", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=93", 2) | Out-Null
$d.Content.Find.Execute("54-3=51", $true, $false, $false, $false, $false, $true, 1, $false, "88-87=1", 2) | Out-Null
$d.Content.Find.Execute("81-17=64", $true, $false, $false, $false, $false, $true, 1, $false, "11+33=44", 2) | Out-Null
$d.Content.Find.Execute("47+39=86", $true, $false, $false, $false, $false, $true, 1, $false, "69-23=46", 2) | Out-Null
$d.Content.Find.Execute("29+7=36", $true, $false, $false, $false, $false, $true, 1, $false, "85+1=86", 2) | Out-Null
$d.Content.Find.Execute("42-37=5", $true, $false, $false, $false, $false, $true, 1, $false, "9+19=28", 2) | Out-Null
$d.Content.Find.Execute("37-28=9", $true, $false, $false, $false, $false, $true, 1, $false, "60-19=41", 2) | Out-Null
$d.Content.Find.Execute("81-33=48", $true, $false, $false, $false, $false, $true, 1, $false, "23-5=18", 2) | Out-Null
$d.Content.Find.Execute("21+31=52", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=94", 2) | Out-Null
$d.Content.Find.Execute("38+13=51", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=53", 2) | Out-Null
$d.Content.Find.Execute("59+17=76", $true, $false, $false, $false, $false, $true, 1, $false, "21+4=25", 2) | Out-Null
$d.Content.Find.Execute("56-18=38", $true, $false, $false, $false, $false, $true, 1, $false, "42-13=29", 2) | Out-Null
$d.Content.Find.Execute("68-57=11", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=4", 2) | Out-Null
$d.Content.Find.Execute("17+55=72", $true, $false, $false, $false, $false, $true, 1, $false, "29+24=53", 2) | Out-Null
$d.Content.Find.Execute("25+10=35", $true, $false, $false, $false, $false, $true, 1, $false, "83-37=46", 2) | Out-Null
$d.Content.Find.Execute("69-37=32", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=37", 2) | Out-Null
$d.Content.Find.Execute("10+25=35", $true, $false, $false, $false, $false, $true, 1, $false, "87-61=26", 2) | Out-Null
$d.Content.Find.Execute("98-85=13", $true, $false, $false, $false, $false, $true, 1, $false, "33+57=90", 2) | Out-Null
$d.Content.Find.Execute("26+37=63", $true, $false, $false, $false, $false, $true, 1, $false, "56-12=44", 2) | Out-Null
$d.Content.Find.Execute("53+16=69", $true, $false, $false, $false, $false, $true, 1, $false, "95-56=39", 2) | Out-Null
$d.Content.Find.Execute("94-62=32", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=2", 2) | Out-Null
$d.Content.Find.Execute("38-24=14", $true, $false, $false, $false, $false, $true, 1, $false, "15+75=90", 2) | Out-Null
$d.Content.Find.Execute("3+25=28", $true, $false, $false, $false, $false, $true, 1, $false, "17+30=47", 2) | Out-Null
$d.Content.Find.Execute("71-32=39", $true, $false, $false, $false, $false, $true, 1, $false, "30-9=21", 2) | Out-Null
$d.Content.Find.Execute("76-23=53", $true, $false, $false, $false, $false, $true, 1, $false, "92-40=52", 2) | Out-Null
$d.Content.Find.Execute("61-26=35", $true, $false, $false, $false, $false, $true, 1, $false, "94-91=3", 2) | Out-Null
$d.Content.Find.Execute("29-4=25", $true, $false, $false, $false, $false, $true, 1, $false, "12+25=37", 2) | Out-Null
$d.Content.Find.Execute("52-22=30", $true, $false, $false, $false, $false, $true, 1, $false, "11+9=20", 2) | Out-Null
$d.Content.Find.Execute("0+24=24", $true, $false, $false, $false, $false, $true, 1, $false, "79-26=53", 2) | Out-Null
$d.Content.Find.Execute("96-64=32", $true, $false, $false, $false, $false, $true, 1, $false, "16+23=39", 2) | Out-Null
$d.Content.Find.Execute("81-19=62", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=3", 2) | Out-Null
$d.Content.Find.Execute("30+47=77", $true, $false, $false, $false, $false, $true, 1, $false, "2+16=18", 2) | Out-Null
$d.Content.Find.Execute("53-31=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+28=47", 2) | Out-Null
$d.Content.Find.Execute("61-17=44", $true, $false, $false, $false, $false, $true, 1, $false, "82+12=94", 2) | Out-Null
$d.Content.Find.Execute("60+0=60", $true, $false, $false, $false, $false, $true, 1, $false, "99-48=51", 2) | Out-Null
$d.Content.Find.Execute("20+2=22", $true, $false, $false, $false, $false, $true, 1, $false, "91-41=50", 2) | Out-Null
$d.Content.Find.Execute("33+34=67", $true, $false, $false, $false, $false, $true, 1, $false, "0+30=30", 2) | Out-Null
$d.Content.Find.Execute("35+37=72", $true, $false, $false, $false, $false, $true, 1, $false, "3+54=57", 2) | Out-Null
$d.Content.Find.Execute("79-5=74", $true, $false, $false, $false, $false, $true, 1, $false, "27-15=12", 2) | Out-Null
$d.Content.Find.Execute("73+2=75", $true, $false, $false, $false, $false, $true, 1, $false, "56-51=5", 2) | Out-Null
$d.Content.Find.Execute("98-56=42", $true, $false, $false, $false, $false, $true, 1, $false, "39+56=95", 2) | Out-Null
$d.Content.Find.Execute("16-6=10", $true, $false, $false, $false, $false, $true, 1, $false, "51+40=91", 2) | Out-Null
$d.Content.Find.Execute("74-42=32", $true, $false, $false, $false, $false, $true, 1, $false, "18+16=34", 2) | Out-Null
$d.Content.Find.Execute("3+69=72", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=27", 2) | Out-Null
$d.Content.Find.Execute("62-33=29", $true, $false, $false, $false, $false, $true, 1, $false, "29+46=75", 2) | Out-Null
$d.Content.Find.Execute("33+17=50", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=22", 2) | Out-Null
$d.Content.Find.Execute("69-32=37", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=64", 2) | Out-Null
$d.Content.Find.Execute("20+22=42", $true, $false, $false, $false, $false, $true, 1, $false, "70+23=93", 2) | Out-Null
$d.Content.Find.Execute("25+13=38", $true, $false, $false, $false, $false, $true, 1, $false, "83-69=14", 2) | Out-Null
$d.Content.Find.Execute("82-21=61", $true, $false, $false, $false, $false, $true, 1, $false, "53-2=51", 2) | Out-Null
$d.Content.Find.Execute("29+16=45", $true, $false, $false, $false, $false, $true, 1, $false, "60+25=85", 2) | Out-Null
$d.Content.Find.Execute("11+81=92", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=7", 2) | Out-Null
$d.Content.Find.Execute("4+65=69", $true, $false, $false, $false, $false, $true, 1, $false, "98-27=71", 2) | Out-Null
$d.Content.Find.Execute("47-34=13", $true, $false, $false, $false, $false, $true, 1, $false, "24+17=41", 2) | Out-Null
$d.Content.Find.Execute("47+7=54", $true, $false, $false, $false, $false, $true, 1, $false, "86-54=32", 2) | Out-Null
$d.Content.Find.Execute("31+39=70", $true, $false, $false, $false, $false, $true, 1, $false, "44-15=29", 2) | Out-Null
$d.Content.Find.Execute("6+43=49", $true, $false, $false, $false, $false, $true, 1, $false, "65+32=97", 2) | Out-Null
$d.Content.Find.Execute("64-58=6", $true, $false, $false, $false, $false, $true, 1, $false, "32+56=88", 2) | Out-Null
$d.Content.Find.Execute("81-38=43", $true, $false, $false, $false, $false, $true, 1, $false, "77-35=42", 2) | Out-Null
$d.Content.Find.Execute("56-47=9", $true, $false, $false, $false, $false, $true, 1, $false, "92-36=56", 2) | Out-Null
$d.Content.Find.Execute("82-19=63", $true, $false, $false, $false, $false, $true, 1, $false, "73+1=74", 2) | Out-Null
$d.Content.Find.Execute("87-27=60", $true, $false, $false, $false, $false, $true, 1, $false, "18+2=20", 2) | Out-Null
$d.Content.Find.Execute("65-12=53", $true, $false, $false, $false, $false, $true, 1, $false, "54-44=10", 2) | Out-Null
$d.Content.Find.Execute("93-63=30", $true, $false, $false, $false, $false, $true, 1, $false, "35+59=94", 2) | Out-Null
$d.Content.Find.Execute("11+43=54", $true, $false, $false, $false, $false, $true, 1, $false, "61-34=27", 2) | Out-Null
$d.Content.Find.Execute("1+96=97", $true, $false, $false, $false, $false, $true, 1, $false, "94-84=10", 2) | Out-Null
$d.Content.Find.Execute("54+3=57", $true, $false, $false, $false, $false, $true, 1, $false, "29+16=45", 2) | Out-Null
$d.Content.Find.Execute("21+29=50", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=21", 2) | Out-Null
$d.Content.Find.Execute("77-68=9", $true, $false, $false, $false, $false, $true, 1, $false, "73-6=67", 2) | Out-Null
$d.Content.Find.Execute("45-37=8", $true, $false, $false, $false, $false, $true, 1, $false, "40+53=93", 2) | Out-Null
$d.Content.Find.Execute("86-60=26", $true, $false, $false, $false, $false, $true, 1, $false, "81-57=24", 2) | Out-Null
$d.Content.Find.Execute("77-18=59", $true, $false, $false, $false, $false, $true, 1, $false, "15+17=32", 2) | Out-Null
$d.Content.Find.Execute("72-65=7", $true, $false, $false, $false, $false, $true, 1, $false, "80-74=6", 2) | Out-Null
$d.Content.Find.Execute("1+35=36", $true, $false, $false, $false, $false, $true, 1, $false, "79-8=71", 2) | Out-Null
$d.Content.Find.Execute("40+35=75", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=95", 2) | Out-Null
$d.Content.Find.Execute("50+35=85", $true, $false, $false, $false, $false, $true, 1, $false, "24-8=16", 2) | Out-Null
$d.Content.Find.Execute("52-42=10", $true, $false, $false, $false, $false, $true, 1, $false, "10+48=58", 2) | Out-Null
$d.Content.Find.Execute("30+5=35", $true, $false, $false, $false, $false, $true, 1, $false, "30-8=22", 2) | Out-Null
$d.Content.Find.Execute("78-44=34", $true, $false, $false, $false, $false, $true, 1, $false, "78-73=5", 2) | Out-Null
$d.Content.Find.Execute("54-11=43", $true, $false, $false, $false, $false, $true, 1, $false, "13+63=76", 2) | Out-Null
$d.Content.Find.Execute("96-70=26", $true, $false, $false, $false, $false, $true, 1, $false, "21+51=72", 2) | Out-Null
$d.Content.Find.Execute("94-65=29", $true, $false, $false, $false, $false, $true, 1, $false, "30+42=72", 2) | Out-Null
$d.Content.Find.Execute("66+22=88", $true, $false, $false, $false, $false, $true, 1, $false, "97-76=21", 2) | Out-Null
$d.Content.Find.Execute("92-30=62", $true, $false, $false, $false, $false, $true, 1, $false, "24-11=13", 2) | Out-Null
$d.Content.Find.Execute("85-84=1", $true, $false, $false, $false, $false, $true, 1, $false, "1+19=20", 2) | Out-Null
$d.Content.Find.Execute("69+21=90", $true, $false, $false, $false, $false, $true, 1, $false, "52+47=99", 2) | Out-Null
$d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "77-63=14", 2) | Out-Null
$d.Content.Find.Execute("84-64=20", $true, $false, $false, $false, $false, $true, 1, $false, "24-11=13", 2) | Out-Null
$d.Content.Find.Execute("83-14=69", $true, $false, $false, $false, $false, $true, 1, $false, "89-81=8", 2) | Out-Null
$d.Content.Find.Execute("17+45=62", $true, $false, $false, $false, $false, $true, 1, $false, "80-33=47", 2) | Out-Null
$d.Content.Find.Execute("67+16=83", $true, $false, $false, $false, $false, $true, 1, $false, "78-13=65", 2) | Out-Null
